$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nIAx_lifeTables3")

$ws.Range("C2").Value = 0.3229350910479403
$ws.Range("D2").Value = 17.75730670488243
$ws.Range("E2").Value = 2906.415819431462
$ws.Range("F2").Value = 17757.30670488243
$ws.Range("G2").Value = 20663.72252431389

$ws.Range("C3").Value = 0.3892355655465126
$ws.Range("D3").Value = 17.85303236555318
$ws.Range("E3").Value = 3503.120089918613
$ws.Range("F3").Value = 17853.03236555318
$ws.Range("G3").Value = 21356.15245547179

$ws.Range("C4").Value = 0.4680540851181295
$ws.Range("D4").Value = 17.24422295904563
$ws.Range("E4").Value = 4212.486766063165
$ws.Range("F4").Value = 17244.22295904563
$ws.Range("G4").Value = 21456.70972510879

$ws.Range("C5").Value = 0.5577562201235239
$ws.Range("D5").Value = 15.80743156200336
$ws.Range("E5").Value = 5019.805981111715
$ws.Range("F5").Value = 15807.43156200336
$ws.Range("G5").Value = 20827.23754311507

$ws.Range("C6").Value = 0.657497716893974
$ws.Range("D6").Value = 13.41054662005103
$ws.Range("E6").Value = 5917.479452045765
$ws.Range("F6").Value = 13410.54662005103
$ws.Range("G6").Value = 19328.0260720968

$ws.Range("C7").Value = 0.7645739323165114
$ws.Range("D7").Value = 10.0038426878967
$ws.Range("E7").Value = 6881.165390848602
$ws.Range("F7").Value = 10003.8426878967
$ws.Range("G7").Value = 16885.0080787453

$ws.Range("C8").Value = 0.8601638560118121
$ws.Range("D8").Value = 6.349693883909262
$ws.Range("E8").Value = 7741.474704106309
$ws.Range("F8").Value = 6349.693883909263
$ws.Range("G8").Value = 14091.16858801557

$ws.Range("C9").Value = 0.2633030297853862
$ws.Range("D9").Value = 16.98193743353636
$ws.Range("E9").Value = 2369.727268068475
$ws.Range("F9").Value = 16981.93743353636
$ws.Range("G9").Value = 19351.66470160483

$ws.Range("C10").Value = 0.3175755584262259
$ws.Range("D10").Value = 17.5888826256501
$ws.Range("E10").Value = 2858.180025836033
$ws.Range("F10").Value = 17588.8826256501
$ws.Range("G10").Value = 20447.06265148614

$ws.Range("C11").Value = 0.3814809155589957
$ws.Range("D11").Value = 17.74837748130697
$ws.Range("E11").Value = 3433.328240030961
$ws.Range("F11").Value = 17748.37748130697
$ws.Range("G11").Value = 21181.70572133793

$ws.Range("C12").Value = 0.4566191823947142
$ws.Range("D12").Value = 17.27748516108446
$ws.Range("E12").Value = 4109.572641552428
$ws.Range("F12").Value = 17277.48516108446
$ws.Range("G12").Value = 21387.05780263689

$ws.Range("C13").Value = 0.5441774788199917
$ws.Range("D13").Value = 15.96664512310385
$ws.Range("E13").Value = 4897.597309379925
$ws.Range("F13").Value = 15966.64512310385
$ws.Range("G13").Value = 20864.24243248377

$ws.Range("C14").Value = 0.6418742833936536
$ws.Range("D14").Value = 13.71151524467597
$ws.Range("E14").Value = 5776.868550542882
$ws.Range("F14").Value = 13711.51524467597
$ws.Range("G14").Value = 19488.38379521885

$ws.Range("C15").Value = 0.7415273408620626
$ws.Range("D15").Value = 10.70491429678136
$ws.Range("E15").Value = 6673.746067758563
$ws.Range("F15").Value = 10704.91429678136
$ws.Range("G15").Value = 17378.66036453992

$ws.Range("C16").Value = 0.3189334484269624
$ws.Range("D16").Value = 17.37641096669558
$ws.Range("E16").Value = 2870.401035842662
$ws.Range("F16").Value = 17376.41096669558
$ws.Range("G16").Value = 20246.81200253824

$ws.Range("C17").Value = 0.3795401182208366
$ws.Range("D17").Value = 17.58701423651813
$ws.Range("E17").Value = 3415.861063987529
$ws.Range("F17").Value = 17587.01423651813
$ws.Range("G17").Value = 21002.87530050566

$ws.Range("C18").Value = 0.4534159030933587
$ws.Range("D18").Value = 17.13056235352797
$ws.Range("E18").Value = 4080.743127840228
$ws.Range("F18").Value = 17130.56235352797
$ws.Range("G18").Value = 21211.3054813682

$ws.Range("C19").Value = 0.5382172297018603
$ws.Range("D19").Value = 15.89521140828897
$ws.Range("E19").Value = 4843.955067316742
$ws.Range("F19").Value = 15895.21140828897
$ws.Range("G19").Value = 20739.16647560571

$ws.Range("C20").Value = 0.6287865630241447
$ws.Range("D20").Value = 13.91823987820568
$ws.Range("E20").Value = 5659.079067217303
$ws.Range("F20").Value = 13918.23987820568
$ws.Range("G20").Value = 19577.31894542299

$ws.Range("C21").Value = 0.7203274626576385
$ws.Range("D21").Value = 11.30890600375432
$ws.Range("E21").Value = 6482.947163918746
$ws.Range("F21").Value = 11308.90600375432
$ws.Range("G21").Value = 17791.85316767307

$ws.Range("C22").Value = 0.801355055737256
$ws.Range("D22").Value = 8.547291667163043
$ws.Range("E22").Value = 7212.195501635304
$ws.Range("F22").Value = 8547.291667163043
$ws.Range("G22").Value = 15759.48716879835
